# remove fast RU, crop, plot HTML in plot_along using laptop ggplot2 3.5.1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the newline-separated "facet" labels in column E with literal
# "<br>" separators (and drop the trailing space before each break),
# matching the updated plot_along() HTML output from ggplot2 3.5.1.
$ws.Range("E2:E15").Value = "**Global**:<br>Implemented by<br>All other countries"
$ws.Range("E16:E29").Value = "**High-income**:<br>All other HICs and<br>not some MICs (such as China)"
$ws.Range("E30:E43").Value = "**International**:<br>Some countries (e.g. EU, UK, Brazil)<br>and not others (e.g. U.S., China)"

# Refreshed confidence-interval numbers for the "All" / global row (row 2)
# and the high-income "All" row (row 16) from the re-run model.
$ws.Range("B2").Value = 73.8161485502559
$ws.Range("C2").Value = 72.4542861110033
$ws.Range("D2").Value = 75.1780109895084

$ws.Range("B16").Value = 69.211130206042
$ws.Range("C16").Value = 67.7830711918507
$ws.Range("D16").Value = 70.6391892202333
